# Fix the format of longitude / latitude values: the values were stored
# with a comma as decimal separator ("13,3761809" / "52,5103735"); they
# should use a dot instead ("13.3761809" / "52.5103735").
#
# Column O holds "Longitude" and column P holds "Latitude" for every data
# row (rows 2-139) on the single worksheet of this workbook, and every row
# shares the exact same two values, so we can fix them all in one shot by
# assigning to the whole column ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 139

$longitudeRange = $ws.Range("O2:O" + $lastRow)
$latitudeRange  = $ws.Range("P2:P" + $lastRow)

$longitudeRange.Value = "13.3761809"
$latitudeRange.Value  = "52.5103735"

# Reflect the view-state changes that were present in the target workbook
# (the cursor/selection ended up on P4 with the sheet scrolled so column F
# is the left-most visible column).
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("P4").Select() | Out-Null
